$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Source data rows (2-18) were reshuffled in the update: each row keeps its
# identity columns (Mercado, Producto, Categoria, Variedad, Calidad, Unidad, Kg/unidad)
# but is now paired with a different record's Fecha / Volumen / Precios / Origen / Precio-Kg.
# The explicit post-edit values below were taken directly from the target diff.

# Row 2 (now holds the record previously in row 13)
$ws.Range("D2").Value = Get-Date -Year 2022 -Month 1 -Day 26 -Hour 0 -Minute 0 -Second 0
$ws.Range("M2").Value = 165
$ws.Range("P2").Value = 6742
$ws.Range("R2").Value = "Provincia de Linares"
$ws.Range("S2").Value = 3371

# Row 3 (now holds the record previously in row 12)
$ws.Range("D3").Value = Get-Date -Year 2021 -Month 1 -Day 15 -Hour 0 -Minute 0 -Second 0
$ws.Range("M3").Value = 45
$ws.Range("R3").Value = "Provincia de Curicó"

# Row 4 (now holds the record previously in row 15)
$ws.Range("D4").Value = Get-Date -Year 2022 -Month 3 -Day 7 -Hour 0 -Minute 0 -Second 0
$ws.Range("M4").Value = 45
$ws.Range("N4").Value = 6000
$ws.Range("O4").Value = 6000
$ws.Range("P4").Value = 6000
$ws.Range("R4").Value = "Provincia de Linares"
$ws.Range("S4").Value = 3000

# Row 5 (now holds the record previously in row 4)
$ws.Range("D5").Value = Get-Date -Year 2022 -Month 1 -Day 25 -Hour 0 -Minute 0 -Second 0
$ws.Range("M5").Value = 80
$ws.Range("N5").Value = 7000
$ws.Range("O5").Value = 7000
$ws.Range("P5").Value = 7000
$ws.Range("R5").Value = "Provincia de Curicó"
$ws.Range("S5").Value = 3500

# Row 6 (now holds the record previously in row 18)
$ws.Range("D6").Value = Get-Date -Year 2022 -Month 1 -Day 28 -Hour 0 -Minute 0 -Second 0
$ws.Range("M6").Value = 60
$ws.Range("N6").Value = 6000
$ws.Range("O6").Value = 6000
$ws.Range("P6").Value = 6000
$ws.Range("S6").Value = 3000

# Row 7 (now holds the record previously in row 16)
$ws.Range("D7").Value = Get-Date -Year 2022 -Month 2 -Day 14 -Hour 0 -Minute 0 -Second 0
$ws.Range("M7").Value = 45
$ws.Range("R7").Value = "Provincia de Linares"

# Row 8 (now holds the record previously in row 10)
$ws.Range("D8").Value = Get-Date -Year 2022 -Month 2 -Day 22 -Hour 0 -Minute 0 -Second 0
$ws.Range("M8").Value = 45
$ws.Range("O8").Value = 6000
$ws.Range("P8").Value = 6000
$ws.Range("R8").Value = "Provincia de Linares"
$ws.Range("S8").Value = 3000

# Row 9 (now holds the record previously in row 3)
$ws.Range("D9").Value = Get-Date -Year 2021 -Month 1 -Day 18 -Hour 0 -Minute 0 -Second 0
$ws.Range("M9").Value = 48
$ws.Range("N9").Value = 6000
$ws.Range("O9").Value = 6000
$ws.Range("P9").Value = 6000
$ws.Range("R9").Value = "Provincia de Linares"
$ws.Range("S9").Value = 3000

# Row 10 (now holds the record previously in row 8)
$ws.Range("D10").Value = Get-Date -Year 2022 -Month 1 -Day 21 -Hour 0 -Minute 0 -Second 0
$ws.Range("M10").Value = 150
$ws.Range("O10").Value = 6500
$ws.Range("P10").Value = 6233
$ws.Range("R10").Value = "Provincia de Curicó"
$ws.Range("S10").Value = 3116

# Row 12 (now holds the record previously in row 9)
$ws.Range("D12").Value = Get-Date -Year 2022 -Month 1 -Day 24 -Hour 0 -Minute 0 -Second 0
$ws.Range("M12").Value = 160
$ws.Range("N12").Value = 6500
$ws.Range("O12").Value = 7000
$ws.Range("P12").Value = 6750
$ws.Range("S12").Value = 3375

# Row 13 (now holds the record previously in row 14)
$ws.Range("D13").Value = Get-Date -Year 2022 -Month 3 -Day 8 -Hour 0 -Minute 0 -Second 0
$ws.Range("M13").Value = 40
$ws.Range("N13").Value = 6000
$ws.Range("O13").Value = 6000
$ws.Range("P13").Value = 6000
$ws.Range("S13").Value = 3000

# Row 14 (now holds the record previously in row 7)
$ws.Range("D14").Value = Get-Date -Year 2023 -Month 2 -Day 2 -Hour 0 -Minute 0 -Second 0
$ws.Range("N14").Value = 7000
$ws.Range("O14").Value = 7000
$ws.Range("P14").Value = 7000
$ws.Range("R14").Value = "Provincia de Curicó"
$ws.Range("S14").Value = 3500

# Row 15 (now holds the record previously in row 17)
$ws.Range("D15").Value = Get-Date -Year 2023 -Month 2 -Day 3 -Hour 0 -Minute 0 -Second 0
$ws.Range("M15").Value = 40
$ws.Range("N15").Value = 7000
$ws.Range("O15").Value = 7000
$ws.Range("P15").Value = 7000
$ws.Range("R15").Value = "Provincia de Curicó"
$ws.Range("S15").Value = 3500

# Row 16 (now holds the record previously in row 5)
$ws.Range("D16").Value = Get-Date -Year 2022 -Month 1 -Day 31 -Hour 0 -Minute 0 -Second 0
$ws.Range("M16").Value = 30
$ws.Range("N16").Value = 8000
$ws.Range("O16").Value = 8000
$ws.Range("P16").Value = 8000
$ws.Range("S16").Value = 4000

# Row 17 (now holds the record previously in row 2)
$ws.Range("D17").Value = Get-Date -Year 2022 -Month 1 -Day 27 -Hour 0 -Minute 0 -Second 0
$ws.Range("M17").Value = 160
$ws.Range("N17").Value = 6500
$ws.Range("P17").Value = 6750
$ws.Range("S17").Value = 3375

# Row 18 (now holds the record previously in row 6)
$ws.Range("D18").Value = Get-Date -Year 2023 -Month 2 -Day 17 -Hour 0 -Minute 0 -Second 0
$ws.Range("M18").Value = 130
$ws.Range("N18").Value = 7000
$ws.Range("O18").Value = 7500
$ws.Range("P18").Value = 7269
$ws.Range("S18").Value = 3634

